$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "Lara"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 10
